$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add the new "Correlation" sheet after the last existing sheet
# ("Weekly Data"), matching the sheet order in the target workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Correlation"

# ---------------------------------------------------------------------
# Headers
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "Box Office Gross"
$ws.Range("D1").Value = "Correlation:"

# ---------------------------------------------------------------------
# Data - mirrors "Weekly Data" rows 5:19, columns C (Rank) and D
# (Box Office Gross), rebased to rows 2:16 here.
# ---------------------------------------------------------------------
$ranks = @(1,1,1,1,2,2,2,3,3,3,3,3,4,4,6)
$gross = @(95412007,108296294,55124155,27811019,15640371,11189286,5874379,3853298,2428830,891660,334824,288499,171178,119527,46421)

for ($i = 0; $i -lt $ranks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ranks[$i]
    $ws.Cells.Item($row, 2).Value = $gross[$i]
}

# ---------------------------------------------------------------------
# Correlation formula
# ---------------------------------------------------------------------
$ws.Range("D2").Formula = "=CORREL(A2:A16,B2:B16)"

# ---------------------------------------------------------------------
# Column B width (matches bestFit width used for the same header /
# data elsewhere in the workbook)
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 13.998697916666666

# ---------------------------------------------------------------------
# Selection / view state for the new sheet
# ---------------------------------------------------------------------
$null = $ws.Activate()
$null = $ws.Range("D3").Select()

# ---------------------------------------------------------------------
# Selection / view state for "Weekly Data" (no longer the active tab)
# ---------------------------------------------------------------------
$wd = $wb.Worksheets.Item("Weekly Data")
$null = $wd.Activate()
$null = $wd.Range("C5:D19").Select()

# Re-activate the Correlation sheet so it ends up as the active tab.
$null = $ws.Activate()
